# Update cover sheet to MPS 2019.3
#
# 1) Split the subtitle run "MPS 2019.2.x + " on the title slide's cover
#    sheet into a new leading run "MPS 2019.3." followed by the
#    (trimmed) remainder "x + ", leaving the rest of the paragraph intact.
# 2) Refresh the cached "datetimeFigureOut" date placeholder text
#    (slide master + every slide layout) from "4-11-2019" to "2-2-2021".

$p = $ppt.ActivePresentation

# --- 1) Cover sheet subtitle -------------------------------------------------

$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Subtitle 4" -and $sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $prefix = $tr.Characters(1, 11)
        if ($prefix.Text -eq "MPS 2019.2.") {
            $prefix.Text = "MPS 2019.3."
        }
    }
}

# --- 2) Date placeholders (slide master + all slide layouts) ---------------

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        if ($shape.HasTextFrame -and $shape.Name -like "*Date Placeholder*") {
            $dtr = $shape.TextFrame.TextRange
            if ($dtr.Text -eq "4-11-2019") {
                $dtr.Text = "2-2-2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
